# Edit script: weekly update for "Hortaliza, Terminal Hortofrutícola Agro Chillán - Repollo"
# A new weekly record is inserted at row 42; all subsequent records (old rows 42-155)
# shift down by one row, and the record that falls off the end becomes the new last
# row (156).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the variable columns (D, I, J, K, L, M, O, P) for rows 42..156.
# Row 42 holds the brand-new record; rows 43..156 hold what used to be in the row above.
$rowsData = @(
    [PSCustomObject]@{ Row=42; D=44525; I='Primera'; J=200; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=43; D=44473; I='Primera'; J=120; K=500; L=550; M=525; O='Región del Maule'; P=525 },
    [PSCustomObject]@{ Row=44; D=44433; I='Primera'; J=300; K=500; L=550; M=525; O='Provincia de Diguillín'; P=525 },
    [PSCustomObject]@{ Row=45; D=44512; I='Primera'; J=300; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=46; D=44460; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=47; D=44438; I='Primera'; J=300; K=600; L=650; M=625; O='Región del Maule'; P=625 },
    [PSCustomObject]@{ Row=48; D=44519; I='Primera'; J=200; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=49; D=44392; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=50; D=44489; I='Primera'; J=200; K=600; L=700; M=650; O='Región del Maule'; P=650 },
    [PSCustomObject]@{ Row=51; D=44434; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=52; D=44449; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=53; D=44399; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=54; D=44298; I='Primera'; J=120; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=55; D=44482; I='Primera'; J=300; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=56; D=44405; I='Primera'; J=160; K=550; L=600; M=575; O='Provincia de Diguillín'; P=575 },
    [PSCustomObject]@{ Row=57; D=44218; I='Primera'; J=2600; K=800; L=900; M=846; O='Región del Maule'; P=846 },
    [PSCustomObject]@{ Row=58; D=44273; I='Primera'; J=3300; K=950; L=1000; M=977; O='Región del Maule'; P=977 },
    [PSCustomObject]@{ Row=59; D=44386; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=60; D=44435; I='Primera'; J=1200; K=500; L=650; M=600; O='Provincia de Diguillín'; P=600 },
    [PSCustomObject]@{ Row=61; D=44328; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=62; D=44442; I='Primera'; J=300; K=650; L=700; M=675; O='Región del Maule'; P=675 },
    [PSCustomObject]@{ Row=63; D=44516; I='Primera'; J=300; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=64; D=44175; I='Primera'; J=120; K=750; L=800; M=775; O='Provincia de Diguillín'; P=775 },
    [PSCustomObject]@{ Row=65; D=44168; I='Primera'; J=120; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=66; D=44203; I='Primera'; J=2800; K=700; L=800; M=754; O='Provincia de Diguillín'; P=754 },
    [PSCustomObject]@{ Row=67; D=44475; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=68; D=44200; I='Primera'; J=2600; K=700; L=800; M=746; O='Provincia de Diguillín'; P=746 },
    [PSCustomObject]@{ Row=69; D=44419; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=70; D=44162; I='Primera'; J=2800; K=900; L=1000; M=946; O='Región de Coquimbo'; P=946 },
    [PSCustomObject]@{ Row=71; D=44357; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=72; D=44202; I='Primera'; J=3200; K=700; L=800; M=753; O='Provincia de Diguillín'; P=753 },
    [PSCustomObject]@{ Row=73; D=44390; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=74; D=44174; I='Primera'; J=120; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=75; D=44293; I='Primera'; J=300; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=76; D=44496; I='Primera'; J=200; K=600; L=700; M=650; O='Región del Maule'; P=650 },
    [PSCustomObject]@{ Row=77; D=44326; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=78; D=44302; I='Primera'; J=300; K=750; L=800; M=775; O='Provincia de Diguillín'; P=775 },
    [PSCustomObject]@{ Row=79; D=44308; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=80; D=44498; I='Primera'; J=300; K=650; L=700; M=675; O='Provincia de Diguillín'; P=675 },
    [PSCustomObject]@{ Row=81; D=44420; I='Primera'; J=120; K=600; L=650; M=625; O='Región del Maule'; P=625 },
    [PSCustomObject]@{ Row=82; D=44398; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=83; D=44396; I='Primera'; J=240; K=500; L=600; M=550; O='Provincia de Diguillín'; P=550 },
    [PSCustomObject]@{ Row=84; D=44321; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=85; D=44208; I='Primera'; J=2900; K=700; L=800; M=748; O='Región del Maule'; P=748 },
    [PSCustomObject]@{ Row=86; D=44349; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=87; D=44477; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=88; D=44487; I='Primera'; J=100; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=89; D=44452; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=90; D=44505; I='Primera'; J=400; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=91; D=44204; I='Primera'; J=2800; K=750; L=800; M=773; O='Provincia de Diguillín'; P=773 },
    [PSCustomObject]@{ Row=92; D=44306; I='Primera'; J=300; K=750; L=800; M=775; O='Provincia de Diguillín'; P=775 },
    [PSCustomObject]@{ Row=93; D=44509; I='Primera'; J=400; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=94; D=44454; I='Primera'; J=300; K=500; L=600; M=550; O='Provincia de Diguillín'; P=550 },
    [PSCustomObject]@{ Row=95; D=44278; I='Primera'; J=300; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=96; D=44265; I='Primera'; J=2700; K=1100; L=1200; M=1156; O='Región Metropolitana'; P=1156 },
    [PSCustomObject]@{ Row=97; D=44494; I='Primera'; J=200; K=600; L=700; M=650; O='Región del Maule'; P=650 },
    [PSCustomObject]@{ Row=98; D=44300; I='Primera'; J=160; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=99; D=44356; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=100; D=44469; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=101; D=44446; I='Primera'; J=300; K=600; L=650; M=625; O='Región del Maule'; P=625 },
    [PSCustomObject]@{ Row=102; D=44323; I='Primera'; J=400; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=103; D=44417; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=104; D=44342; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=105; D=44523; I='Primera'; J=240; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=106; D=44406; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=107; D=44295; I='Primera'; J=120; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=108; D=44270; I='Primera'; J=2700; K=1000; L=1100; M=1048; O='Región del Maule'; P=1048 },
    [PSCustomObject]@{ Row=109; D=44363; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=110; D=44299; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=111; D=44372; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=112; D=44372; I='Segunda'; J=80; K=500; L=500; M=500; O='Provincia de Diguillín'; P=500 },
    [PSCustomObject]@{ Row=113; D=44403; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=114; D=44169; I='Primera'; J=120; K=750; L=800; M=775; O='Provincia de Diguillín'; P=775 },
    [PSCustomObject]@{ Row=115; D=44195; I='Primera'; J=600; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=116; D=44376; I='Primera'; J=120; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=117; D=44524; I='Primera'; J=240; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=118; D=44172; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=119; D=44421; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=120; D=44426; I='Primera'; J=300; K=500; L=550; M=525; O='Provincia de Diguillín'; P=525 },
    [PSCustomObject]@{ Row=121; D=44448; I='Primera'; J=300; K=600; L=650; M=625; O='Región del Maule'; P=625 },
    [PSCustomObject]@{ Row=122; D=44362; I='Primera'; J=120; K=500; L=550; M=525; O='Provincia de Diguillín'; P=525 },
    [PSCustomObject]@{ Row=123; D=44210; I='Primera'; J=2600; K=700; L=800; M=750; O='Provincia de Diguillín'; P=750 },
    [PSCustomObject]@{ Row=124; D=44176; I='Primera'; J=300; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=125; D=44301; I='Primera'; J=300; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=126; D=44407; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=127; D=44284; I='Primera'; J=120; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=128; D=44441; I='Primera'; J=300; K=600; L=650; M=625; O='Región del Maule'; P=625 },
    [PSCustomObject]@{ Row=129; D=44504; I='Primera'; J=360; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=130; D=44350; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=131; D=44382; I='Primera'; J=160; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=132; D=44329; I='Primera'; J=300; K=500; L=550; M=525; O='Provincia de Diguillín'; P=525 },
    [PSCustomObject]@{ Row=133; D=44522; I='Primera'; J=300; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=134; D=44491; I='Primera'; J=200; K=600; L=700; M=650; O='Región del Maule'; P=650 },
    [PSCustomObject]@{ Row=135; D=44305; I='Primera'; J=120; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=136; D=44225; I='Primera'; J=2800; K=900; L=1000; M=946; O='Región del Maule'; P=946 },
    [PSCustomObject]@{ Row=137; D=44447; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=138; D=44425; I='Primera'; J=300; K=500; L=550; M=525; O='Provincia de Diguillín'; P=525 },
    [PSCustomObject]@{ Row=139; D=44322; I='Primera'; J=600; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=140; D=44495; I='Primera'; J=180; K=600; L=700; M=650; O='Región del Maule'; P=650 },
    [PSCustomObject]@{ Row=141; D=44232; I='Primera'; J=300; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=142; D=44327; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=143; D=44510; I='Primera'; J=360; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=144; D=44161; I='Primera'; J=2600; K=950; L=1000; M=977; O='Región del Maule'; P=977 },
    [PSCustomObject]@{ Row=145; D=44468; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=146; D=44517; I='Primera'; J=200; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=147; D=44391; I='Primera'; J=160; K=500; L=600; M=550; O='Provincia de Diguillín'; P=550 },
    [PSCustomObject]@{ Row=148; D=44236; I='Primera'; J=300; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=149; D=44340; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=150; D=44515; I='Primera'; J=240; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 },
    [PSCustomObject]@{ Row=151; D=44330; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=152; D=44432; I='Primera'; J=300; K=600; L=650; M=625; O='Provincia de Diguillín'; P=625 },
    [PSCustomObject]@{ Row=153; D=44181; I='Primera'; J=160; K=800; L=850; M=825; O='Provincia de Diguillín'; P=825 },
    [PSCustomObject]@{ Row=154; D=44194; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=155; D=44307; I='Primera'; J=160; K=700; L=750; M=725; O='Provincia de Diguillín'; P=725 },
    [PSCustomObject]@{ Row=156; D=44508; I='Primera'; J=400; K=600; L=700; M=650; O='Provincia de Diguillín'; P=650 }
)

foreach ($r in $rowsData) {
    $ws.Cells.Item($r.Row, 4).Value2  = $r.D   # Fecha
    $ws.Cells.Item($r.Row, 9).Value2  = $r.I   # Calidad
    $ws.Cells.Item($r.Row, 10).Value2 = $r.J   # Volumen
    $ws.Cells.Item($r.Row, 11).Value2 = $r.K   # Precio minimo
    $ws.Cells.Item($r.Row, 12).Value2 = $r.L   # Precio maximo
    $ws.Cells.Item($r.Row, 13).Value2 = $r.M   # Precio promedio ponderado
    $ws.Cells.Item($r.Row, 15).Value2 = $r.O   # Origen
    $ws.Cells.Item($r.Row, 16).Value2 = $r.P   # Precio $/Kg
}

# Row 156 is a brand-new row; fill in the columns that stay constant across every
# record of this subset, copied from the row directly above it.
$ws.Cells.Item(156, 1).Value2  = $ws.Cells.Item(155, 1).Value2    # Mercado ID
$ws.Cells.Item(156, 2).Value2  = $ws.Cells.Item(155, 2).Value2    # Mercado
$ws.Cells.Item(156, 3).Value2  = $ws.Cells.Item(155, 3).Value2    # Region
$ws.Cells.Item(156, 5).Value2  = $ws.Cells.Item(155, 5).Value2    # Codreg
$ws.Cells.Item(156, 6).Value2  = $ws.Cells.Item(155, 6).Value2    # Categoria ID
$ws.Cells.Item(156, 7).Value2  = $ws.Cells.Item(155, 7).Value2    # Categoria
$ws.Cells.Item(156, 8).Value2  = $ws.Cells.Item(155, 8).Value2    # Variedad
$ws.Cells.Item(156, 14).Value2 = $ws.Cells.Item(155, 14).Value2   # Unidad de comercializacion
$ws.Cells.Item(156, 17).Value2 = $ws.Cells.Item(155, 17).Value2   # Kg o Unidades
$ws.Cells.Item(156, 18).Value2 = $ws.Cells.Item(155, 18).Value2   # Clasificacion

# Match the date number format used by the rest of column D.
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(155, 4).NumberFormat
